$d = $word.ActiveDocument

# Phase 1: replace each original value with a unique placeholder token
# to avoid any cross-matching between old/new values across cells.
$d.Content.Find.Execute("2024-05-20 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH000@@", 2) | Out-Null
$d.Content.Find.Execute("72-6=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH001@@", 2) | Out-Null
$d.Content.Find.Execute("26+41=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH002@@", 2) | Out-Null
$d.Content.Find.Execute("50+29=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH003@@", 2) | Out-Null
$d.Content.Find.Execute("27+11=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH004@@", 2) | Out-Null
$d.Content.Find.Execute("75+3=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH005@@", 2) | Out-Null
$d.Content.Find.Execute("83-36=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH006@@", 2) | Out-Null
$d.Content.Find.Execute("84-55=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH007@@", 2) | Out-Null
$d.Content.Find.Execute("73-10=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH008@@", 2) | Out-Null
$d.Content.Find.Execute("3+72=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH009@@", 2) | Out-Null
$d.Content.Find.Execute("66+8=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH010@@", 2) | Out-Null
$d.Content.Find.Execute("44-5=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH011@@", 2) | Out-Null
$d.Content.Find.Execute("10+44=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH012@@", 2) | Out-Null
$d.Content.Find.Execute("4+54=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH013@@", 2) | Out-Null
$d.Content.Find.Execute("41+8=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH014@@", 2) | Out-Null
$d.Content.Find.Execute("4+42=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH015@@", 2) | Out-Null
$d.Content.Find.Execute("19-12=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH016@@", 2) | Out-Null
$d.Content.Find.Execute("22+45=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH017@@", 2) | Out-Null
$d.Content.Find.Execute("77+14=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH018@@", 2) | Out-Null
$d.Content.Find.Execute("21-7=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH019@@", 2) | Out-Null
$d.Content.Find.Execute("72-37=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH020@@", 2) | Out-Null
$d.Content.Find.Execute("37+4=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH021@@", 2) | Out-Null
$d.Content.Find.Execute("46+50=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH022@@", 2) | Out-Null
$d.Content.Find.Execute("17+16=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH023@@", 2) | Out-Null
$d.Content.Find.Execute("33-0=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH024@@", 2) | Out-Null
$d.Content.Find.Execute("78-70=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH025@@", 2) | Out-Null
$d.Content.Find.Execute("71-28=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH026@@", 2) | Out-Null
$d.Content.Find.Execute("23+1=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH027@@", 2) | Out-Null
$d.Content.Find.Execute("71+3=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH028@@", 2) | Out-Null
$d.Content.Find.Execute("68-55=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH029@@", 2) | Out-Null
$d.Content.Find.Execute("28+64=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH030@@", 2) | Out-Null
$d.Content.Find.Execute("31+52=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH031@@", 2) | Out-Null
$d.Content.Find.Execute("91-38=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH032@@", 2) | Out-Null
$d.Content.Find.Execute("4+82=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH033@@", 2) | Out-Null
$d.Content.Find.Execute("25-12=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH034@@", 2) | Out-Null
$d.Content.Find.Execute("3+42=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH035@@", 2) | Out-Null
$d.Content.Find.Execute("3+68=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH036@@", 2) | Out-Null
$d.Content.Find.Execute("42-36=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH037@@", 2) | Out-Null
$d.Content.Find.Execute("39-19=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH038@@", 2) | Out-Null
$d.Content.Find.Execute("84-34=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH039@@", 2) | Out-Null
$d.Content.Find.Execute("55-43=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH040@@", 2) | Out-Null
$d.Content.Find.Execute("55+16=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH041@@", 2) | Out-Null
$d.Content.Find.Execute("36-28=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH042@@", 2) | Out-Null
$d.Content.Find.Execute("84-82=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH043@@", 2) | Out-Null
$d.Content.Find.Execute("59+3=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH044@@", 2) | Out-Null
$d.Content.Find.Execute("37+58=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH045@@", 2) | Out-Null
$d.Content.Find.Execute("54-10=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH046@@", 2) | Out-Null
$d.Content.Find.Execute("22-19=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH047@@", 2) | Out-Null
$d.Content.Find.Execute("59-59=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH048@@", 2) | Out-Null
$d.Content.Find.Execute("81-35=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH049@@", 2) | Out-Null
$d.Content.Find.Execute("64-25=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH050@@", 2) | Out-Null
$d.Content.Find.Execute("25+44=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH051@@", 2) | Out-Null
$d.Content.Find.Execute("50-33=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH052@@", 2) | Out-Null
$d.Content.Find.Execute("19+67=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH053@@", 2) | Out-Null
$d.Content.Find.Execute("39+3=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH054@@", 2) | Out-Null
$d.Content.Find.Execute("61-14=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH055@@", 2) | Out-Null
$d.Content.Find.Execute("74-20=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH056@@", 2) | Out-Null
$d.Content.Find.Execute("91-8=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH057@@", 2) | Out-Null
$d.Content.Find.Execute("56+10=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH058@@", 2) | Out-Null
$d.Content.Find.Execute("23-1=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH059@@", 2) | Out-Null
$d.Content.Find.Execute("15+38=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH060@@", 2) | Out-Null
$d.Content.Find.Execute("89-40=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH061@@", 2) | Out-Null
$d.Content.Find.Execute("41+36=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH062@@", 2) | Out-Null
$d.Content.Find.Execute("24+36=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH063@@", 2) | Out-Null
$d.Content.Find.Execute("90-22=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH064@@", 2) | Out-Null
$d.Content.Find.Execute("50-47=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH065@@", 2) | Out-Null
$d.Content.Find.Execute("4+7=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH066@@", 2) | Out-Null
$d.Content.Find.Execute("13+27=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH067@@", 2) | Out-Null
$d.Content.Find.Execute("47+4=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH068@@", 2) | Out-Null
$d.Content.Find.Execute("43-5=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH069@@", 2) | Out-Null
$d.Content.Find.Execute("97-73=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH070@@", 2) | Out-Null
$d.Content.Find.Execute("64+18=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH071@@", 2) | Out-Null
$d.Content.Find.Execute("57-12=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH072@@", 2) | Out-Null
$d.Content.Find.Execute("27-24=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH073@@", 2) | Out-Null
$d.Content.Find.Execute("66+9=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH074@@", 2) | Out-Null
$d.Content.Find.Execute("77-57=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH075@@", 2) | Out-Null
$d.Content.Find.Execute("80-62=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH076@@", 2) | Out-Null
$d.Content.Find.Execute("26+43=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH077@@", 2) | Out-Null
$d.Content.Find.Execute("7+17=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH078@@", 2) | Out-Null
$d.Content.Find.Execute("96-92=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH079@@", 2) | Out-Null
$d.Content.Find.Execute("79-23=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH080@@", 2) | Out-Null
$d.Content.Find.Execute("78-39=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH081@@", 2) | Out-Null
$d.Content.Find.Execute("48+49=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH082@@", 2) | Out-Null
$d.Content.Find.Execute("31+42=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH083@@", 2) | Out-Null
$d.Content.Find.Execute("26+20=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH084@@", 2) | Out-Null
$d.Content.Find.Execute("1+82=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH085@@", 2) | Out-Null
$d.Content.Find.Execute("90-34=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH086@@", 2) | Out-Null
$d.Content.Find.Execute("0+31=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH087@@", 2) | Out-Null
$d.Content.Find.Execute("90-18=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH088@@", 2) | Out-Null
$d.Content.Find.Execute("16+79=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH089@@", 2) | Out-Null
$d.Content.Find.Execute("16+38=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH090@@", 2) | Out-Null
$d.Content.Find.Execute("55-54=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH091@@", 2) | Out-Null
$d.Content.Find.Execute("20+72=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH092@@", 2) | Out-Null
$d.Content.Find.Execute("82+4=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH093@@", 2) | Out-Null
$d.Content.Find.Execute("87-18=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH094@@", 2) | Out-Null
$d.Content.Find.Execute("89+1=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH095@@", 2) | Out-Null
$d.Content.Find.Execute("29+65=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH096@@", 2) | Out-Null
$d.Content.Find.Execute("50+44=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH097@@", 2) | Out-Null
$d.Content.Find.Execute("97-17=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH098@@", 2) | Out-Null
$d.Content.Find.Execute("90-72=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH099@@", 2) | Out-Null
$d.Content.Find.Execute("94-22=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH100@@", 2) | Out-Null

# Phase 2: replace each placeholder with the final new value
$d.Content.Find.Execute("@@PH000@@", $true, $false, $false, $false, $false, $true, 1, $false, "2024-05-21 Tuesday", 2) | Out-Null
$d.Content.Find.Execute("@@PH001@@", $true, $false, $false, $false, $false, $true, 1, $false, "29-2=", 2) | Out-Null
$d.Content.Find.Execute("@@PH002@@", $true, $false, $false, $false, $false, $true, 1, $false, "2+66=", 2) | Out-Null
$d.Content.Find.Execute("@@PH003@@", $true, $false, $false, $false, $false, $true, 1, $false, "46-1=", 2) | Out-Null
$d.Content.Find.Execute("@@PH004@@", $true, $false, $false, $false, $false, $true, 1, $false, "93-40=", 2) | Out-Null
$d.Content.Find.Execute("@@PH005@@", $true, $false, $false, $false, $false, $true, 1, $false, "61-44=", 2) | Out-Null
$d.Content.Find.Execute("@@PH006@@", $true, $false, $false, $false, $false, $true, 1, $false, "80-20=", 2) | Out-Null
$d.Content.Find.Execute("@@PH007@@", $true, $false, $false, $false, $false, $true, 1, $false, "38+27=", 2) | Out-Null
$d.Content.Find.Execute("@@PH008@@", $true, $false, $false, $false, $false, $true, 1, $false, "52-17=", 2) | Out-Null
$d.Content.Find.Execute("@@PH009@@", $true, $false, $false, $false, $false, $true, 1, $false, "72+5=", 2) | Out-Null
$d.Content.Find.Execute("@@PH010@@", $true, $false, $false, $false, $false, $true, 1, $false, "8+60=", 2) | Out-Null
$d.Content.Find.Execute("@@PH011@@", $true, $false, $false, $false, $false, $true, 1, $false, "38+47=", 2) | Out-Null
$d.Content.Find.Execute("@@PH012@@", $true, $false, $false, $false, $false, $true, 1, $false, "27+5=", 2) | Out-Null
$d.Content.Find.Execute("@@PH013@@", $true, $false, $false, $false, $false, $true, 1, $false, "39-4=", 2) | Out-Null
$d.Content.Find.Execute("@@PH014@@", $true, $false, $false, $false, $false, $true, 1, $false, "17+59=", 2) | Out-Null
$d.Content.Find.Execute("@@PH015@@", $true, $false, $false, $false, $false, $true, 1, $false, "76-11=", 2) | Out-Null
$d.Content.Find.Execute("@@PH016@@", $true, $false, $false, $false, $false, $true, 1, $false, "50-23=", 2) | Out-Null
$d.Content.Find.Execute("@@PH017@@", $true, $false, $false, $false, $false, $true, 1, $false, "48+51=", 2) | Out-Null
$d.Content.Find.Execute("@@PH018@@", $true, $false, $false, $false, $false, $true, 1, $false, "38-12=", 2) | Out-Null
$d.Content.Find.Execute("@@PH019@@", $true, $false, $false, $false, $false, $true, 1, $false, "66+20=", 2) | Out-Null
$d.Content.Find.Execute("@@PH020@@", $true, $false, $false, $false, $false, $true, 1, $false, "74+18=", 2) | Out-Null
$d.Content.Find.Execute("@@PH021@@", $true, $false, $false, $false, $false, $true, 1, $false, "11+2=", 2) | Out-Null
$d.Content.Find.Execute("@@PH022@@", $true, $false, $false, $false, $false, $true, 1, $false, "60-57=", 2) | Out-Null
$d.Content.Find.Execute("@@PH023@@", $true, $false, $false, $false, $false, $true, 1, $false, "83-24=", 2) | Out-Null
$d.Content.Find.Execute("@@PH024@@", $true, $false, $false, $false, $false, $true, 1, $false, "30+16=", 2) | Out-Null
$d.Content.Find.Execute("@@PH025@@", $true, $false, $false, $false, $false, $true, 1, $false, "10+33=", 2) | Out-Null
$d.Content.Find.Execute("@@PH026@@", $true, $false, $false, $false, $false, $true, 1, $false, "42+15=", 2) | Out-Null
$d.Content.Find.Execute("@@PH027@@", $true, $false, $false, $false, $false, $true, 1, $false, "13+10=", 2) | Out-Null
$d.Content.Find.Execute("@@PH028@@", $true, $false, $false, $false, $false, $true, 1, $false, "59+33=", 2) | Out-Null
$d.Content.Find.Execute("@@PH029@@", $true, $false, $false, $false, $false, $true, 1, $false, "75-74=", 2) | Out-Null
$d.Content.Find.Execute("@@PH030@@", $true, $false, $false, $false, $false, $true, 1, $false, "61-42=", 2) | Out-Null
$d.Content.Find.Execute("@@PH031@@", $true, $false, $false, $false, $false, $true, 1, $false, "19-7=", 2) | Out-Null
$d.Content.Find.Execute("@@PH032@@", $true, $false, $false, $false, $false, $true, 1, $false, "23+68=", 2) | Out-Null
$d.Content.Find.Execute("@@PH033@@", $true, $false, $false, $false, $false, $true, 1, $false, "60+33=", 2) | Out-Null
$d.Content.Find.Execute("@@PH034@@", $true, $false, $false, $false, $false, $true, 1, $false, "53-17=", 2) | Out-Null
$d.Content.Find.Execute("@@PH035@@", $true, $false, $false, $false, $false, $true, 1, $false, "79-73=", 2) | Out-Null
$d.Content.Find.Execute("@@PH036@@", $true, $false, $false, $false, $false, $true, 1, $false, "99-17=", 2) | Out-Null
$d.Content.Find.Execute("@@PH037@@", $true, $false, $false, $false, $false, $true, 1, $false, "97-58=", 2) | Out-Null
$d.Content.Find.Execute("@@PH038@@", $true, $false, $false, $false, $false, $true, 1, $false, "88-56=", 2) | Out-Null
$d.Content.Find.Execute("@@PH039@@", $true, $false, $false, $false, $false, $true, 1, $false, "94-49=", 2) | Out-Null
$d.Content.Find.Execute("@@PH040@@", $true, $false, $false, $false, $false, $true, 1, $false, "33+36=", 2) | Out-Null
$d.Content.Find.Execute("@@PH041@@", $true, $false, $false, $false, $false, $true, 1, $false, "67+4=", 2) | Out-Null
$d.Content.Find.Execute("@@PH042@@", $true, $false, $false, $false, $false, $true, 1, $false, "62-21=", 2) | Out-Null
$d.Content.Find.Execute("@@PH043@@", $true, $false, $false, $false, $false, $true, 1, $false, "51+16=", 2) | Out-Null
$d.Content.Find.Execute("@@PH044@@", $true, $false, $false, $false, $false, $true, 1, $false, "97-91=", 2) | Out-Null
$d.Content.Find.Execute("@@PH045@@", $true, $false, $false, $false, $false, $true, 1, $false, "0+61=", 2) | Out-Null
$d.Content.Find.Execute("@@PH046@@", $true, $false, $false, $false, $false, $true, 1, $false, "94-28=", 2) | Out-Null
$d.Content.Find.Execute("@@PH047@@", $true, $false, $false, $false, $false, $true, 1, $false, "61-30=", 2) | Out-Null
$d.Content.Find.Execute("@@PH048@@", $true, $false, $false, $false, $false, $true, 1, $false, "84-72=", 2) | Out-Null
$d.Content.Find.Execute("@@PH049@@", $true, $false, $false, $false, $false, $true, 1, $false, "56+14=", 2) | Out-Null
$d.Content.Find.Execute("@@PH050@@", $true, $false, $false, $false, $false, $true, 1, $false, "87-80=", 2) | Out-Null
$d.Content.Find.Execute("@@PH051@@", $true, $false, $false, $false, $false, $true, 1, $false, "67-65=", 2) | Out-Null
$d.Content.Find.Execute("@@PH052@@", $true, $false, $false, $false, $false, $true, 1, $false, "11+35=", 2) | Out-Null
$d.Content.Find.Execute("@@PH053@@", $true, $false, $false, $false, $false, $true, 1, $false, "1+29=", 2) | Out-Null
$d.Content.Find.Execute("@@PH054@@", $true, $false, $false, $false, $false, $true, 1, $false, "24-6=", 2) | Out-Null
$d.Content.Find.Execute("@@PH055@@", $true, $false, $false, $false, $false, $true, 1, $false, "72-63=", 2) | Out-Null
$d.Content.Find.Execute("@@PH056@@", $true, $false, $false, $false, $false, $true, 1, $false, "63-59=", 2) | Out-Null
$d.Content.Find.Execute("@@PH057@@", $true, $false, $false, $false, $false, $true, 1, $false, "60-44=", 2) | Out-Null
$d.Content.Find.Execute("@@PH058@@", $true, $false, $false, $false, $false, $true, 1, $false, "0+50=", 2) | Out-Null
$d.Content.Find.Execute("@@PH059@@", $true, $false, $false, $false, $false, $true, 1, $false, "26+33=", 2) | Out-Null
$d.Content.Find.Execute("@@PH060@@", $true, $false, $false, $false, $false, $true, 1, $false, "36+37=", 2) | Out-Null
$d.Content.Find.Execute("@@PH061@@", $true, $false, $false, $false, $false, $true, 1, $false, "93-32=", 2) | Out-Null
$d.Content.Find.Execute("@@PH062@@", $true, $false, $false, $false, $false, $true, 1, $false, "94-49=", 2) | Out-Null
$d.Content.Find.Execute("@@PH063@@", $true, $false, $false, $false, $false, $true, 1, $false, "11+71=", 2) | Out-Null
$d.Content.Find.Execute("@@PH064@@", $true, $false, $false, $false, $false, $true, 1, $false, "0+47=", 2) | Out-Null
$d.Content.Find.Execute("@@PH065@@", $true, $false, $false, $false, $false, $true, 1, $false, "44+44=", 2) | Out-Null
$d.Content.Find.Execute("@@PH066@@", $true, $false, $false, $false, $false, $true, 1, $false, "14+15=", 2) | Out-Null
$d.Content.Find.Execute("@@PH067@@", $true, $false, $false, $false, $false, $true, 1, $false, "61-36=", 2) | Out-Null
$d.Content.Find.Execute("@@PH068@@", $true, $false, $false, $false, $false, $true, 1, $false, "22+73=", 2) | Out-Null
$d.Content.Find.Execute("@@PH069@@", $true, $false, $false, $false, $false, $true, 1, $false, "34+33=", 2) | Out-Null
$d.Content.Find.Execute("@@PH070@@", $true, $false, $false, $false, $false, $true, 1, $false, "5+42=", 2) | Out-Null
$d.Content.Find.Execute("@@PH071@@", $true, $false, $false, $false, $false, $true, 1, $false, "68-52=", 2) | Out-Null
$d.Content.Find.Execute("@@PH072@@", $true, $false, $false, $false, $false, $true, 1, $false, "46-9=", 2) | Out-Null
$d.Content.Find.Execute("@@PH073@@", $true, $false, $false, $false, $false, $true, 1, $false, "32+9=", 2) | Out-Null
$d.Content.Find.Execute("@@PH074@@", $true, $false, $false, $false, $false, $true, 1, $false, "79-3=", 2) | Out-Null
$d.Content.Find.Execute("@@PH075@@", $true, $false, $false, $false, $false, $true, 1, $false, "17-8=", 2) | Out-Null
$d.Content.Find.Execute("@@PH076@@", $true, $false, $false, $false, $false, $true, 1, $false, "7+24=", 2) | Out-Null
$d.Content.Find.Execute("@@PH077@@", $true, $false, $false, $false, $false, $true, 1, $false, "88-74=", 2) | Out-Null
$d.Content.Find.Execute("@@PH078@@", $true, $false, $false, $false, $false, $true, 1, $false, "38-30=", 2) | Out-Null
$d.Content.Find.Execute("@@PH079@@", $true, $false, $false, $false, $false, $true, 1, $false, "94-14=", 2) | Out-Null
$d.Content.Find.Execute("@@PH080@@", $true, $false, $false, $false, $false, $true, 1, $false, "24+70=", 2) | Out-Null
$d.Content.Find.Execute("@@PH081@@", $true, $false, $false, $false, $false, $true, 1, $false, "47-7=", 2) | Out-Null
$d.Content.Find.Execute("@@PH082@@", $true, $false, $false, $false, $false, $true, 1, $false, "95-64=", 2) | Out-Null
$d.Content.Find.Execute("@@PH083@@", $true, $false, $false, $false, $false, $true, 1, $false, "76+0=", 2) | Out-Null
$d.Content.Find.Execute("@@PH084@@", $true, $false, $false, $false, $false, $true, 1, $false, "60-41=", 2) | Out-Null
$d.Content.Find.Execute("@@PH085@@", $true, $false, $false, $false, $false, $true, 1, $false, "3+23=", 2) | Out-Null
$d.Content.Find.Execute("@@PH086@@", $true, $false, $false, $false, $false, $true, 1, $false, "20+9=", 2) | Out-Null
$d.Content.Find.Execute("@@PH087@@", $true, $false, $false, $false, $false, $true, 1, $false, "73-56=", 2) | Out-Null
$d.Content.Find.Execute("@@PH088@@", $true, $false, $false, $false, $false, $true, 1, $false, "79-32=", 2) | Out-Null
$d.Content.Find.Execute("@@PH089@@", $true, $false, $false, $false, $false, $true, 1, $false, "28+42=", 2) | Out-Null
$d.Content.Find.Execute("@@PH090@@", $true, $false, $false, $false, $false, $true, 1, $false, "91-43=", 2) | Out-Null
$d.Content.Find.Execute("@@PH091@@", $true, $false, $false, $false, $false, $true, 1, $false, "29-7=", 2) | Out-Null
$d.Content.Find.Execute("@@PH092@@", $true, $false, $false, $false, $false, $true, 1, $false, "21+59=", 2) | Out-Null
$d.Content.Find.Execute("@@PH093@@", $true, $false, $false, $false, $false, $true, 1, $false, "52-5=", 2) | Out-Null
$d.Content.Find.Execute("@@PH094@@", $true, $false, $false, $false, $false, $true, 1, $false, "23+31=", 2) | Out-Null
$d.Content.Find.Execute("@@PH095@@", $true, $false, $false, $false, $false, $true, 1, $false, "1+61=", 2) | Out-Null
$d.Content.Find.Execute("@@PH096@@", $true, $false, $false, $false, $false, $true, 1, $false, "18+48=", 2) | Out-Null
$d.Content.Find.Execute("@@PH097@@", $true, $false, $false, $false, $false, $true, 1, $false, "73-54=", 2) | Out-Null
$d.Content.Find.Execute("@@PH098@@", $true, $false, $false, $false, $false, $true, 1, $false, "30+37=", 2) | Out-Null
$d.Content.Find.Execute("@@PH099@@", $true, $false, $false, $false, $false, $true, 1, $false, "77-50=", 2) | Out-Null
$d.Content.Find.Execute("@@PH100@@", $true, $false, $false, $false, $false, $true, 1, $false, "19+42=", 2) | Out-Null
